$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B67 was stored as a text "3"; convert it to a real number 3
$ws.Range("B67").Value = 3

# Append new row 68 with the new annotation entry
$ws.Range("A68").Value = "Ruilin"

# B68 must stay textual ("2"), not be auto-coerced to a number
$ws.Range("B68").NumberFormat = "@"
$ws.Range("B68").Value = "2"
$ws.Range("B68").Style = "Normal"

$ws.Range("C68").Value = "really bad"
$ws.Range("D68").Value = "CRT"
$ws.Range("E68").Value = "OTH"
$ws.Range("F68").Value = "3222e19c-371b-4610-a09f-eba8d4490b26"
$ws.Range("G68").Value = "rJTGkKxAZ_annotated.xlsx"
$ws.Range("H68").Value = "This section is really bad."
